$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. "Here is our score..." paragraph: pseudo-likelihood -> likelihood
#    conditioned on w.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Here is our score.  The first term is a pseudo-likelihood.  The second term follows the SparsityBoost framework (explain that).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Here is our score.  The first term is the likelihood conditioned on w.  The second term follows the SparsityBoost framework (explain that).",
    2)

# ------------------------------------------------------------------
# 2. Heading "Pseudo-likelihood term" -> "Likelihood term"
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Pseudo-likelihood term",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Likelihood term",
    2)

# ------------------------------------------------------------------
# 3. "linear-in-RKHS assumption, plus Gaussian noise" -> remove the
#    redundant second "assumption"
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Equivalent to ML under the assumption linear-in-RKHS assumption, plus Gaussian noise.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Equivalent to ML under the assumption linear-in-RKHS, plus Gaussian noise.",
    2)

# ------------------------------------------------------------------
# 4. Move the "_GoBack" bookmark: remove it from between "Present
#    theorem! " and "Show a plot." and re-add it at the very end of
#    the document (end of the CONCLUSION paragraph), after adding a
#    new trailing sentence there.
# ------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# Append the new sentence to the end of the CONCLUSION paragraph
# (i.e. right before the very last paragraph mark of the document).
$full = $d.Content.Text
$endIdx = $full.Length - 1
$insPoint = $d.Range($endIdx, $endIdx)
$insPoint.InsertBefore("Can compute a BDE-like score, i.e. marginalize over prior over w.  ")

# Re-add the "_GoBack" bookmark immediately after the new sentence.
# Placing a collapsed range exactly on the paragraph-mark boundary is
# unreliable, so temporarily insert a one-character marker after the
# target position, anchor the bookmark just before it, then remove
# the marker again - the bookmark stays correctly anchored.
$full2 = $d.Content.Text
$pos = $full2.Length - 1
$marker = $d.Range($pos, $pos)
$marker.InsertBefore("X")

$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$markerRange = $d.Range($pos, $pos + 1)
$markerRange.Text = ""
